# ============================================================================
# Commit: [ADDITIONAL SCRAPING] added scraping code for extra browling
# attributes and excel sheets
#
# 1) Clean up now-redundant empty placeholder cells left on "ODI Batting
#    Extra" (columns B-E were written as empty strings by the old scraper
#    for rows that had no batting-extra data; the new scraper only emits a
#    cell when it actually has a value).
# 2) Add a brand-new "ODI Bowling Extra" sheet (mirrors "ODI Batting Extra"
#    but scraped from the "ODI Bowling" sheet) with MATCH_CODE,
#    MAIDEN_OVERS and PERCENT_WICKETS_OF_ALL columns.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Step 1: strip the stray empty cells out of "ODI Batting Extra"
# ----------------------------------------------------------------------
$wsBattingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$wsBattingExtra.Range("B2:E2").ClearContents()
$wsBattingExtra.Range("B8:E8").ClearContents()
$wsBattingExtra.Range("E11:E11").ClearContents()
$wsBattingExtra.Range("B15:E15").ClearContents()
$wsBattingExtra.Range("B19:E19").ClearContents()
$wsBattingExtra.Range("B20:E20").ClearContents()
$wsBattingExtra.Range("C21:E21").ClearContents()

# ----------------------------------------------------------------------
# Step 2: add the new "ODI Bowling Extra" sheet after "ODI Batting Extra"
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsBowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsBowlingExtra.Name = "ODI Bowling Extra"

# Header row - values first, then copy the header formatting (bold, border,
# centered) from the sibling "ODI Batting Extra" sheet so the new sheet
# matches the workbook's existing header style.
$wsBowlingExtra.Range("A1").Value = "MATCH_CODE"
$wsBowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$wsBowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$wsBattingExtra.Range("A1:C1").Copy()
$wsBowlingExtra.Range("A1:C1").PasteSpecial(-4122)

# Data rows: MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL, matching
# the tail of the "ODI Bowling" sheet (last 20 matches). All three columns
# are text, so the range is pre-formatted as Text before values are
# assigned (otherwise "0"/"4447"/"10.00%" would be auto-parsed as numbers
# or percentages) and the style is reset back to Normal afterwards so the
# data rows keep the workbook's default (unstyled) look.
$data = @(
    @("4447", "0", ""),
    @("4463", "0", "10.00%"),
    @("4464", "", ""),
    @("4465", "0", ""),
    @("4477", "3", "50.00%"),
    @("4479", "0", "20.00%"),
    @("4481", "0", "10.00%"),
    @("4537", "1", "20.00%"),
    @("4538", "", ""),
    @("4539", "0", "10.00%"),
    @("4550", "0", ""),
    @("4557", "2", "10.00%"),
    @("4559", "0", "20.00%"),
    @("4679", "2", "50.00%"),
    @("4682", "1", "20.00%"),
    @("4685", "", ""),
    @("4711", "0", "10.00%"),
    @("4713", "0", "10.00%"),
    @("4717", "0", "40.00%"),
    @("4726", "", "")
)

$rowNum = 2
foreach ($rowData in $data) {
    $matchCode = $rowData[0]
    $maidenOvers = $rowData[1]
    $percentWickets = $rowData[2]

    $rowRange = $wsBowlingExtra.Range("A" + $rowNum + ":C" + $rowNum)
    $rowRange.NumberFormat = "@"

    $wsBowlingExtra.Range("A" + $rowNum).Value = $matchCode
    if ($maidenOvers -ne "") {
        $wsBowlingExtra.Range("B" + $rowNum).Value = $maidenOvers
    }
    if ($percentWickets -ne "") {
        $wsBowlingExtra.Range("C" + $rowNum).Value = $percentWickets
    }

    $rowRange.Style = "Normal"

    $rowNum = $rowNum + 1
}

Write-Output "Done"
